# Update column F (dSF) values for specific rows to reflect a repulled
# dataset. Only the "dSF" column changes; all other columns are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F
$updates = @{
    2  = 0
    7  = -1
    8  = 3
    14 = 1
    16 = 2
    24 = 6
    26 = -2
    29 = -3
    35 = 1
    36 = 6
    42 = -2
    52 = 0
    58 = 2
    63 = -5
    67 = 4
    73 = -2
    76 = -5
    78 = 3
    85 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
